$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheet and tables (Migraatio -> Migraatiot, Table1 -> Migraatiot, Table2 -> Tilat)
$ws1.Name = "Migraatiot"
$lo = $ws1.ListObjects.Item(1)
$lo.Name = "Migraatiot"
$lo2 = $ws2.ListObjects.Item(1)
$lo2.Name = "Tilat"

# Grow the Migraatiot table from 9 data rows to 12 data rows (10 -> 13 incl. header)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Replace the old placeholder/test row with fresh test-fixture data
$ws1.Range("A2").Value = "1.2.246.578.5.1.2293640800.1682339657"
$ws1.Range("B2").Value = "JATKOPAATOS2"
$ws1.Range("C2").Value = "TESTI-ASIA-123"
$ws1.Range("D2").Value = 44611
$ws1.Range("E2").Value = "TESTI-ASIA-243"
$ws1.Range("F2").Value = 44606

# Clear the now-unused oid/Tila/asianumero cells on the rest of the rows
$ws1.Range("A3:C13").Clear()
$ws1.Range("E3:E13").Clear()

# Switch the date columns to an unambiguous ISO format
$ws1.Range("D2:D13").NumberFormat = "yyyy-mm-dd;@"
$ws1.Range("F2:F13").NumberFormat = "yyyy-mm-dd;@"

# Add the list-based data validation for Tila, driven off the Tilat table
$rng = $ws1.Range("B2:B13")
$rng.Validation.Add(3, 1, 1, '=INDIRECT("Tilat[Tila]")')

# touch J12 (matches the stray used-range cell left behind by the original edit)
$ws1.Range("J12").Value = "x"
$ws1.Range("J12").Value = ""

Write-Host "done"
